# Apply updated cryptocurrency price/volume figures to sheet1 (Coin list)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Plain text values (Excel will not reinterpret these as numbers) ---
$ws.Range("D2").Value = "74.751.21"
$ws.Range("E2").Value = "  +0.78%  "
$ws.Range("D3").Value = "2.802.79"
$ws.Range("E3").Value = "  +6.67%  "
$ws.Range("E4").Value = "  +0.14%  "
$ws.Range("E5").Value = "  -0.01%  "
$ws.Range("E6").Value = "  +1.83%  "
$ws.Range("E7").Value = "  +0.10%  "
$ws.Range("E8").Value = "  +2.52%  "
$ws.Range("E9").Value = "  -5.37%  "
$ws.Range("D10").Value = "2.803.95"
$ws.Range("E11").Value = "  -1.60%  "
$ws.Range("E12").Value = "  +3.40%  "
$ws.Range("E13").Value = "  +2.19%  "
$ws.Range("D14").Value = "3.325.95"
$ws.Range("E14").Value = "  +6.89%  "
$ws.Range("D15").Value = "74.709.16"
$ws.Range("E15").Value = "  +0.96%  "
$ws.Range("E16").Value = "  -1.66%  "
$ws.Range("E17").Value = "  +1.48%  "
$ws.Range("D18").Value = "2.806.01"
$ws.Range("E18").Value = "  +6.16%  "
$ws.Range("E19").Value = "  -0.92%  "
$ws.Range("E20").Value = "  +3.56%  "
$ws.Range("E21").Value = "  +0.87%  "
$ws.Range("E22").Value = "  -2.30%  "
$ws.Range("E24").Value = "  -0.07%  "
$ws.Range("E25").Value = "  +0.93%  "
$ws.Range("E26").Value = "  +5.23%  "
$ws.Range("E27").Value = "  +7.08%  "
$ws.Range("E28").Value = "  -0.57%  "
$ws.Range("E29").Value = "  +9.07%  "
$ws.Range("E30").Value = "  +0.53%  "
$ws.Range("E31").Value = "  -2.05%  "
$ws.Range("E32").Value = "  -0.85%  "
$ws.Range("E33").Value = "  -0.24%  "
$ws.Range("E34").Value = "  +2.09%  "
$ws.Range("E35").Value = "  +0.13%  "
$ws.Range("E36").Value = "  +0.02%  "
$ws.Range("E37").Value = "  +3.65%  "
$ws.Range("E38").Value = "  -1.25%  "
$ws.Range("E39").Value = "  +0.37%  "
$ws.Range("E40").Value = "  +15.04%  "
$ws.Range("E41").Value = "  -0.01%  "
$ws.Range("E42").Value = "  +3.98%  "
$ws.Range("E43").Value = "  +0.94%  "
$ws.Range("E44").Value = "  -1.03%  "
$ws.Range("E45").Value = "  +1.98%  "
$ws.Range("E46").Value = "  +2.88%  "
$ws.Range("E47").Value = "  +0.05%  "
$ws.Range("E48").Value = "  -3.04%  "
$ws.Range("E49").Value = "  +8.12%  "
$ws.Range("E50").Value = "  +2.17%  "
$ws.Range("E51").Value = "  +8.14%  "

# --- Numeric-looking text values: must be forced to stay text. ---
# Trick: write a formula returning the literal string, then PasteSpecial
# the value back over itself (xlPasteValues = -4163). This converts the
# formula result into a plain text cell without altering any cell style.
$c = $ws.Range("D5")
$c.Formula = "=""186.30"""
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Range("D6")
$c.Formula = "=""592.62"""
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Range("D8")
$c.Formula = "=""0.547"""
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Range("D9")
$c.Formula = "=""0.191"""
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Range("D11")
$c.Formula = "=""0.160"""
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Range("D13")
$c.Formula = "=""4.87"""
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Range("D16")
$c.Formula = "=""0.0000186"""
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Range("D17")
$c.Formula = "=""26.66"""
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Range("D19")
$c.Formula = "=""8.92"""
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Range("D20")
$c.Formula = "=""12.24"""
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Range("D21")
$c.Formula = "=""375.83"""
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Range("D22")
$c.Formula = "=""2.26"""
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Range("D23")
$c.Formula = "=""4.06"""
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Range("D24")
$c.Formula = "=""1.00"""
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Range("D25")
$c.Formula = "=""70.76"""
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Range("D26")
$c.Formula = "=""9.83"""
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Range("D28")
$c.Formula = "=""4.13"""
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Range("D29")
$c.Formula = "=""0.0000103"""
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Range("D30")
$c.Formula = "=""0.999"""
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Range("D31")
$c.Formula = "=""516.94"""
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Range("D32")
$c.Formula = "=""1.38"""
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Range("D33")
$c.Formula = "=""7.71"""
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Range("D35")
$c.Formula = "=""1.00"""
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Range("D36")
$c.Formula = "=""163.20"""
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Range("D37")
$c.Formula = "=""19.86"""
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Range("D40")
$c.Formula = "=""184.64"""
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Range("D42")
$c.Formula = "=""0.340"""
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Range("D43")
$c.Formula = "=""4.98"""
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Range("D44")
$c.Formula = "=""1.65"""
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Range("D46")
$c.Formula = "=""39.89"""
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Range("D47")
$c.Formula = "=""0.0857"""
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Range("D48")
$c.Formula = "=""2.31"""
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Range("D49")
$c.Formula = "=""0.572"""
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Range("D51")
$c.Formula = "=""0.633"""
$c.Copy()
$c.PasteSpecial(-4163)

$excel.CutCopyMode = 0
